$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.59643625221294
$ws.Range("C2").Value = 12.34594967997777
$ws.Range("D2").Value = 4.111064188861349
$ws.Range("F2").Value = 21.00532472282775
$ws.Range("G2").Value = 23.30210648127286
$ws.Range("H2").Value = 12.75543502850116
$ws.Range("L2").Value = 10.63791345242836
$ws.Range("N2").Value = 16.70572584640052
$ws.Range("O2").Value = 18.70743494873291

$ws.Range("B3").Value = 14.043647788201
$ws.Range("C3").Value = 12.29562328194713
$ws.Range("D3").Value = 4.04837263468885
$ws.Range("F3").Value = 20.97220032541269
$ws.Range("G3").Value = 23.21336587745706
$ws.Range("H3").Value = 12.79130364475449
$ws.Range("L3").Value = 10.60912633936949
$ws.Range("N3").Value = 16.7469185550067
$ws.Range("O3").Value = 18.74170346069273

$ws.Range("B4").Value = 13.69438327828685
$ws.Range("C4").Value = 12.26522658134587
$ws.Range("D4").Value = 4.008874866494646
$ws.Range("F4").Value = 20.95849463701322
$ws.Range("G4").Value = 23.16829677946267
$ws.Range("H4").Value = 12.81575437365076
$ws.Range("L4").Value = 10.59359583528836
$ws.Range("N4").Value = 16.77399418392744
$ws.Range("O4").Value = 18.76773952502286

$ws.Range("B5").Value = 13.5497901020752
$ws.Range("C5").Value = 12.25297309837267
$ws.Range("D5").Value = 3.992536677548201
$ws.Range("F5").Value = 20.95458157079884
$ws.Range("G5").Value = 23.15231375024963
$ws.Range("H5").Value = 12.82632781373001
$ws.Range("L5").Value = 10.58781109709847
$ws.Range("N5").Value = 16.78547679600859
$ws.Range("O5").Value = 18.7796022381818

$ws.Range("B6").Value = 13.52565035091116
$ws.Range("C6").Value = 12.25094666845947
$ws.Range("D6").Value = 3.989809389416263
$ws.Range("F6").Value = 20.95403288905206
$ws.Range("G6").Value = 23.14980405804128
$ws.Range("H6").Value = 12.82812031809306
$ws.Range("L6").Value = 10.58688353346424
$ws.Range("N6").Value = 16.78741061842045
$ws.Range("O6").Value = 18.78164759527327

$ws.Range("B7").Value = 13.69244211857708
$ws.Range("C7").Value = 12.26506077762326
$ws.Range("D7").Value = 4.008655492277345
$ws.Range("F7").Value = 20.958435089468
$ws.Range("G7").Value = 23.16807156165824
$ws.Range("H7").Value = 12.81589450371487
$ws.Range("L7").Value = 10.5935156115376
$ws.Range("N7").Value = 16.77414722324891
$ws.Range("O7").Value = 18.76789444179764

$ws.Range("B8").Value = 14.40799929779707
$ws.Range("C8").Value = 12.3284966081234
$ws.Range("D8").Value = 4.089662692154343
$ws.Range("F8").Value = 20.99252987012304
$ws.Range("G8").Value = 23.26956412994504
$ws.Range("H8").Value = 12.76729829402121
$ws.Range("L8").Value = 10.62754546846141
$ws.Range("N8").Value = 16.71955941297818
$ws.Range("O8").Value = 18.71821219231641

$ws.Range("B9").Value = 15.72487033920883
$ws.Range("C9").Value = 12.45659450220896
$ws.Range("D9").Value = 4.240117445972235
$ws.Range("F9").Value = 21.11174967548959
$ws.Range("G9").Value = 23.54247145499841
$ws.Range("H9").Value = 12.69129400438746
$ws.Range("L9").Value = 10.71106651806844
$ws.Range("N9").Value = 16.6266309211507
$ws.Range("O9").Value = 18.66054473840141

$ws.Range("B10").Value = 16.63027583817763
$ws.Range("C10").Value = 12.55256054940213
$ws.Range("D10").Value = 4.344993546664909
$ws.Range("F10").Value = 21.23078158430192
$ws.Range("G10").Value = 23.78660437111249
$ws.Range("H10").Value = 12.64725868142667
$ws.Range("L10").Value = 10.78232316975521
$ws.Range("N10").Value = 16.56692207113197
$ws.Range("O10").Value = 18.64255771302783

$ws.Range("B11").Value = 17.02707150672954
$ws.Range("C11").Value = 12.59653251094771
$ws.Range("D11").Value = 4.39136826425701
$ws.Range("F11").Value = 21.29161954365954
$ws.Range("G11").Value = 23.90675983750628
$ws.Range("H11").Value = 12.62979751880094
$ws.Range("L11").Value = 10.81680500871575
$ws.Range("N11").Value = 16.54161053427248
$ws.Range("O11").Value = 18.63968890993211

$ws.Range("B12").Value = 17.17504813240705
$ws.Range("C12").Value = 12.61322147623457
$ws.Range("D12").Value = 4.408728972193067
$ws.Range("F12").Value = 21.31560470845686
$ws.Range("G12").Value = 23.9535307575172
$ws.Range("H12").Value = 12.6235556853692
$ws.Range("L12").Value = 10.83015183107262
$ws.Range("N12").Value = 16.53229117501231
$ws.Range("O12").Value = 18.63936731728579

$ws.Range("B13").Value = 17.14328181810799
$ws.Range("C13").Value = 12.60962565496748
$ws.Range("D13").Value = 4.404999077191278
$ws.Range("F13").Value = 21.31039720342238
$ws.Range("G13").Value = 23.94340192312758
$ws.Range("H13").Value = 12.62488349552778
$ws.Range("L13").Value = 10.82726461802301
$ws.Range("N13").Value = 16.53428646229492
$ws.Range("O13").Value = 18.63940256056462

$ws.Range("B14").Value = 17.03929194663018
$ws.Range("C14").Value = 12.59790481058646
$ws.Range("D14").Value = 4.392800603759668
$ws.Range("F14").Value = 21.29357392626423
$ws.Range("G14").Value = 23.91058247026425
$ws.Range("H14").Value = 12.62927657389095
$ws.Range("L14").Value = 10.81789730870902
$ws.Range("N14").Value = 16.54083850537437
$ws.Range("O14").Value = 18.63964712358368

$ws.Range("B15").Value = 16.97529496238866
$ws.Range("C15").Value = 12.59073013777463
$ws.Range("D15").Value = 4.385302343108242
$ws.Range("F15").Value = 21.28339206106615
$ws.Range("G15").Value = 23.89064391310832
$ws.Range("H15").Value = 12.63201570763068
$ws.Range("L15").Value = 10.81219699424912
$ws.Range("N15").Value = 16.54488639140207
$ws.Range("O15").Value = 18.63989652995908

$ws.Range("B16").Value = 16.60403087994295
$ws.Range("C16").Value = 12.54969258399096
$ws.Range("D16").Value = 4.341935264314492
$ws.Range("F16").Value = 21.22693906693506
$ws.Range("G16").Value = 23.77893186126757
$ws.Range("H16").Value = 12.6484515972286
$ws.Range("L16").Value = 10.7801106243763
$ws.Range("N16").Value = 16.56861342934671
$ws.Range("O16").Value = 18.64285218986563

$ws.Range("B17").Value = 16.37232513401253
$ws.Range("C17").Value = 12.52459292808516
$ws.Range("D17").Value = 4.314983067856137
$ws.Range("F17").Value = 21.19401038396319
$ws.Range("G17").Value = 23.71270462764084
$ws.Range("H17").Value = 12.65919340050719
$ws.Range("L17").Value = 10.76095055456489
$ws.Range("N17").Value = 16.58364274635523
$ws.Range("O17").Value = 18.64602699195344

$ws.Range("B18").Value = 16.23764063726595
$ws.Range("C18").Value = 12.51018653214551
$ws.Range("D18").Value = 4.29935580381018
$ws.Range("F18").Value = 21.17570134542888
$ws.Range("G18").Value = 23.67547130414699
$ws.Range("H18").Value = 12.66561373511283
$ws.Range("L18").Value = 10.75012523142013
$ws.Range("N18").Value = 16.59246140173237
$ws.Range("O18").Value = 18.64835319017021

$ws.Range("B19").Value = 16.19179971661744
$ws.Range("C19").Value = 12.50531420765162
$ws.Range("D19").Value = 4.294043467835632
$ws.Range("F19").Value = 21.16961096129877
$ws.Range("G19").Value = 23.6630133402817
$ws.Range("H19").Value = 12.66782908113005
$ws.Range("L19").Value = 10.74649369542389
$ws.Range("N19").Value = 16.59547717978267
$ws.Range("O19").Value = 18.64922666972377

$ws.Range("B20").Value = 16.39713776443723
$ws.Range("C20").Value = 12.52726175343428
$ws.Range("D20").Value = 4.317865183342387
$ws.Range("F20").Value = 21.1974505285827
$ws.Range("G20").Value = 23.71966600675216
$ws.Range("H20").Value = 12.65802487298349
$ws.Range("L20").Value = 10.76297004646192
$ws.Range("N20").Value = 16.58202482560882
$ws.Range("O20").Value = 18.64563725990671

$ws.Range("B21").Value = 17.06989903162643
$ws.Range("C21").Value = 12.60134654475782
$ws.Range("D21").Value = 4.396389095708286
$ws.Range("F21").Value = 21.29848975176329
$ws.Range("G21").Value = 23.92018816956091
$ws.Range("H21").Value = 12.62797616469365
$ws.Range("L21").Value = 10.82064092824552
$ws.Range("N21").Value = 16.53890680792134
$ws.Range("O21").Value = 18.6395545316341

$ws.Range("B22").Value = 17.49624835307193
$ws.Range("C22").Value = 12.64998156090708
$ws.Range("D22").Value = 4.446537325058643
$ws.Range("F22").Value = 21.37003738123851
$ws.Range("G22").Value = 24.058627723144
$ws.Range("H22").Value = 12.61049642187858
$ws.Range("L22").Value = 10.86001446496953
$ws.Range("N22").Value = 16.5122744165287
$ws.Range("O22").Value = 18.64003662629442

$ws.Range("B23").Value = 17.26995170788683
$ws.Range("C23").Value = 12.624006926606
$ws.Range("D23").Value = 4.419882194434817
$ws.Range("F23").Value = 21.33135198938579
$ws.Range("G23").Value = 23.98407727460434
$ws.Range("H23").Value = 12.6196279556723
$ws.Range("L23").Value = 10.83884886274504
$ws.Range("N23").Value = 16.52634717003002
$ws.Range("O23").Value = 18.63937138980362

$ws.Range("B24").Value = 16.38592455595202
$ws.Range("C24").Value = 12.52605510297888
$ws.Range("D24").Value = 4.316562589697749
$ws.Range("F24").Value = 21.19589330029786
$ws.Range("G24").Value = 23.71651614292567
$ws.Range("H24").Value = 12.65855240200381
$ws.Range("L24").Value = 10.76205644158107
$ws.Range("N24").Value = 16.58275573276018
$ws.Range("O24").Value = 18.64581189735037

$ws.Range("B25").Value = 15.37891983264817
$ws.Range("C25").Value = 12.42158568977704
$ws.Range("D25").Value = 4.200372302001971
$ws.Range("F25").Value = 21.07393427587497
$ws.Range("G25").Value = 23.46086379784549
$ws.Range("H25").Value = 12.70978502880779
$ws.Range("L25").Value = 10.6867077095683
$ws.Range("N25").Value = 16.650263412098
$ws.Range("O25").Value = 18.67187198338943
